$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 741.1111
$ws.Cells.Item(28, 9).Value = 679.0714
$ws.Cells.Item(28, 11).Value = 679.0714
$ws.Cells.Item(28, 13).Value = -194.0714

$ws.Cells.Item(41, 8).Value = 302.65
$ws.Cells.Item(41, 9).Value = 287.125
$ws.Cells.Item(41, 11).Value = 287.125
$ws.Cells.Item(41, 13).Value = 152.875

$ws.Cells.Item(62, 8).Value = 6595.0415
$ws.Cells.Item(62, 10).Value = 4237.5
$ws.Cells.Item(62, 12).Value = 4237.5
$ws.Cells.Item(62, 14).Value = -5485.5

$ws.Cells.Item(65, 8).Value = 6595.0415
$ws.Cells.Item(65, 10).Value = 4237.5
$ws.Cells.Item(65, 12).Value = 21187.5
$ws.Cells.Item(65, 14).Value = -27427.5

$ws.Cells.Item(76, 8).Value = 4253.1113
$ws.Cells.Item(76, 9).Value = 4184.875
$ws.Cells.Item(76, 11).Value = 4184.875
$ws.Cells.Item(76, 13).Value = -3869.875

$ws.Cells.Item(79, 8).Value = 4253.1113
$ws.Cells.Item(79, 9).Value = 4184.875
$ws.Cells.Item(79, 11).Value = 4184.875
$ws.Cells.Item(79, 13).Value = -3092.875

$ws.Cells.Item(88, 8).Value = 3800.3333
$ws.Cells.Item(88, 9).Value = 5003
$ws.Cells.Item(88, 10).Value = 3199
$ws.Cells.Item(88, 11).Value = 5003
$ws.Cells.Item(88, 12).Value = 3199
$ws.Cells.Item(88, 13).Value = -4597
$ws.Cells.Item(88, 14).Value = -4011

$ws.Cells.Item(91, 8).Value = 3800.3333
$ws.Cells.Item(91, 9).Value = 5003
$ws.Cells.Item(91, 10).Value = 3199
$ws.Cells.Item(91, 11).Value = 5003
$ws.Cells.Item(91, 12).Value = 3199
$ws.Cells.Item(91, 13).Value = -3599
$ws.Cells.Item(91, 14).Value = -6007

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3905.8071
$ws.Cells.Item(32, 9).Value = 3905.8071
$ws.Cells.Item(32, 11).Value = 3905.8071
$ws.Cells.Item(32, 13).Value = -3618.8071

$ws.Cells.Item(74, 8).Value = 13633.393
$ws.Cells.Item(74, 9).Value = 13359.143
$ws.Cells.Item(74, 11).Value = 13359.143
$ws.Cells.Item(74, 13).Value = -12485.143

$ws.Cells.Item(77, 8).Value = 13633.393
$ws.Cells.Item(77, 9).Value = 13359.143
$ws.Cells.Item(77, 11).Value = 66795.715
$ws.Cells.Item(77, 13).Value = -62427.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 230237.02
$ws.Cells.Item(86, 9).Value = 1001839.7
$ws.Cells.Item(86, 10).Value = 3295.0588
$ws.Cells.Item(86, 11).Value = 1001839.7
$ws.Cells.Item(86, 12).Value = 3295.0588
$ws.Cells.Item(86, 13).Value = -1000716.7
$ws.Cells.Item(86, 14).Value = -5541.0588

$ws.Cells.Item(89, 8).Value = 230237.02
$ws.Cells.Item(89, 9).Value = 1001839.7
$ws.Cells.Item(89, 10).Value = 3295.0588
$ws.Cells.Item(89, 11).Value = 5009198.5
$ws.Cells.Item(89, 12).Value = 16475.294
$ws.Cells.Item(89, 13).Value = -5003582.5
$ws.Cells.Item(89, 14).Value = -27707.294

$ws.Cells.Item(94, 8).Value = 953.2258
$ws.Cells.Item(94, 9).Value = 826.8
$ws.Cells.Item(94, 10).Value = 1480
$ws.Cells.Item(94, 11).Value = 826.8
$ws.Cells.Item(94, 12).Value = 1480
$ws.Cells.Item(94, 13).Value = -375.8
$ws.Cells.Item(94, 14).Value = -2382

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3637.9285
$ws.Cells.Item(31, 9).Value = 3105.1
$ws.Cells.Item(31, 10).Value = 4970
$ws.Cells.Item(31, 11).Value = 3105.1
$ws.Cells.Item(31, 12).Value = 4970
$ws.Cells.Item(31, 13).Value = -2810.1
$ws.Cells.Item(31, 14).Value = -5560

$ws.Cells.Item(34, 8).Value = 3637.9285
$ws.Cells.Item(34, 9).Value = 3105.1
$ws.Cells.Item(34, 10).Value = 4970
$ws.Cells.Item(34, 11).Value = 3105.1
$ws.Cells.Item(34, 12).Value = 4970
$ws.Cells.Item(34, 13).Value = -2903.1
$ws.Cells.Item(34, 14).Value = -5374

$ws.Cells.Item(60, 8).Value = 22401.75
$ws.Cells.Item(60, 9).Value = 1623
$ws.Cells.Item(60, 10).Value = 29328
$ws.Cells.Item(60, 11).Value = 1623
$ws.Cells.Item(60, 12).Value = 29328
$ws.Cells.Item(60, 13).Value = -1112
$ws.Cells.Item(60, 14).Value = -30350

$ws.Cells.Item(94, 8).Value = 1564
$ws.Cells.Item(94, 9).Value = 1564
$ws.Cells.Item(94, 11).Value = 1564
$ws.Cells.Item(94, 13).Value = -1113

$ws.Cells.Item(99, 8).Value = 8623.062
$ws.Cells.Item(99, 10).Value = 10980.774
$ws.Cells.Item(99, 12).Value = 10980.774
$ws.Cells.Item(99, 14).Value = -13976.774

$ws.Cells.Item(126, 8).Value = 8623.062
$ws.Cells.Item(126, 10).Value = 10980.774
$ws.Cells.Item(126, 12).Value = 32942.322
$ws.Cells.Item(126, 14).Value = -37882.322

$ws.Cells.Item(132, 8).Value = 32004.424
$ws.Cells.Item(132, 9).Value = 21722.715
$ws.Cells.Item(132, 11).Value = 65168.145
$ws.Cells.Item(132, 13).Value = -62638.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 3034.6206
$ws.Cells.Item(131, 10).Value = 3909.6843
$ws.Cells.Item(131, 12).Value = 11729.0529
$ws.Cells.Item(131, 14).Value = -21809.0529

$ws.Cells.Item(132, 8).Value = 5884222.5
$ws.Cells.Item(132, 10).Value = 12501860
$ws.Cells.Item(132, 12).Value = 112516740
$ws.Cells.Item(132, 14).Value = -112521800

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 34160.918
$ws.Cells.Item(43, 9).Value = 14999.5
$ws.Cells.Item(43, 10).Value = 37993.2
$ws.Cells.Item(43, 11).Value = 14999.5
$ws.Cells.Item(43, 12).Value = 37993.2
$ws.Cells.Item(43, 13).Value = -14848.5

$ws.Cells.Item(102, 8).Value = 1950.5
$ws.Cells.Item(102, 9).Value = 1950.5
$ws.Cells.Item(102, 11).Value = 1950.5
$ws.Cells.Item(102, 13).Value = -328.5

$ws.Cells.Item(113, 8).Value = 113595.164
$ws.Cells.Item(113, 9).Value = 156299.08
$ws.Cells.Item(113, 11).Value = 156299.08
$ws.Cells.Item(113, 13).Value = -154129.08

$ws.Cells.Item(134, 8).Value = 51661.668
$ws.Cells.Item(134, 10).Value = 51661.668
$ws.Cells.Item(134, 12).Value = 154985.004
$ws.Cells.Item(134, 14).Value = -160055.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 432.2
$ws.Cells.Item(55, 9).Value = 300.33334
$ws.Cells.Item(55, 10).Value = 488.7143
$ws.Cells.Item(55, 11).Value = 300.33334
$ws.Cells.Item(55, 12).Value = 488.7143
$ws.Cells.Item(55, 13).Value = -127.33334
$ws.Cells.Item(55, 14).Value = -834.7143

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 468685.88
$ws.Cells.Item(122, 9).Value = 587287.6
$ws.Cells.Item(122, 10).Value = 17999.2
$ws.Cells.Item(122, 11).Value = 1761862.8
$ws.Cells.Item(122, 12).Value = 53997.60000000001
$ws.Cells.Item(122, 13).Value = -1759412.8
$ws.Cells.Item(122, 14).Value = -58897.60000000001

$ws.Cells.Item(126, 8).Value = 7149.923
$ws.Cells.Item(126, 9).Value = 5245
$ws.Cells.Item(126, 11).Value = 15735
$ws.Cells.Item(126, 13).Value = -13265

$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(141, 8).Value = 62999.332
$ws.Cells.Item(141, 10).Value = 62999.332
$ws.Cells.Item(141, 12).Value = 62999.332
$ws.Cells.Item(141, 14).Value = -73359.33199999999
